$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay as plain text,
# matching the workbooks convention of storing Price as text (t="inlineStr").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.991.83'
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.957.07'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.20'
$ws.Range("E5").Value = '  -1.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4879'
$ws.Range("E7").Value = '  +0.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2958'
$ws.Range("E8").Value = '  +0.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06813'
$ws.Range("E9").Value = '  +0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.16'
$ws.Range("E10").Value = '  -1.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '106.98'
$ws.Range("E11").Value = '  -2.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.954.66'
$ws.Range("E12").Value = '  -0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07828'
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.488'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7029'
$ws.Range("E15").Value = '  +2.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.55'
$ws.Range("E16").Value = '  -3.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.012.72'
$ws.Range("E17").Value = '  +1.11%  '

$ws.Range("B18").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C18").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.249.48'
$ws.Range("E18").Value = '  +0.72%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.18'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007677'
$ws.Range("E20").Value = '  -0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9992'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.502'
$ws.Range("E22").Value = '  -2.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9976'
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.504'
$ws.Range("E24").Value = '  -1.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.773'
$ws.Range("E25").Value = '  -1.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.22'
$ws.Range("E26").Value = '  -0.75%  '

$ws.Range("E27").Value = '  -1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.204'
$ws.Range("E28").Value = '  +0.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1056'
$ws.Range("E29").Value = '  -0.96%  '

$ws.Range("E30").Value = '  -2.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.586'
$ws.Range("E31").Value = '  -1.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.605'
$ws.Range("E32").Value = '  -2.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.434'
$ws.Range("E33").Value = '  -0.14%  '

$ws.Range("E34").Value = '  -3.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7614'
$ws.Range("E35").Value = '  -1.28%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.174'
$ws.Range("E36").Value = '  -0.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.730'
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02013'
$ws.Range("E38").Value = '  -1.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.697'
$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("E40").Value = '  +2.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '77.07'
$ws.Range("E41").Value = '  +9.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.114'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8889'
$ws.Range("E43").Value = '  +1.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4468'
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.04'
$ws.Range("E45").Value = '  -0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.118'
$ws.Range("E46").Value = '  +8.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003.17'
$ws.Range("E48").Value = '  +9.36%  '

$ws.Range("E49").Value = '  -1.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.350'
$ws.Range("E50").Value = '  -0.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.88'
$ws.Range("E51").Value = '  -0.65%  '
